$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of quote data appended to the sheet (row 4).
# Force text format first so Excel doesn't auto-convert the date-like
# string "2025-09-04" into a date serial number.
$row = 4
$ws.Range("A$row" + ":C$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2025-09-04"
$ws.Range("B$row").Value = "21:21:01"
$ws.Range("C$row").Value = "1.00 EUR = 1589.8516 ARS"
